# Rename worksheets from Russian names to English "Worksheet N" names,
# and update the date in the shared report title string.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "Worksheet 1"
$wb.Worksheets.Item(2).Name = "Worksheet 2"
$wb.Worksheets.Item(3).Name = "Worksheet 3"

$ws1 = $wb.Worksheets.Item(1)
$newTitle = "Date: 05-10-2018 - Department: Sales department"
$ws1.Range("A2").Value = $newTitle
$ws1.Range("F4").Value = $newTitle
$ws1.Range("J4").Value = $newTitle
$ws1.Range("F11").Value = $newTitle
$ws1.Range("J11").Value = $newTitle
$ws1.Range("A31").Value = $newTitle
